$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("nationality") is removed; column G ("address") shifts left into its place,
# carrying along its precise custom width. Column H (numbers) shifts left too, so we
# re-insert a blank column to push it back to H.
$ws.Columns.Item(6).Delete() | Out-Null
$ws.Columns.Item(7).Insert() | Out-Null

# Restore the nationality values into the (now empty) column G for rows 2 and 3.
# Row 1's nationality is simply dropped (no G1 cell at all), matching the target.
$ws.Range("G2").Value = "rumana"
$ws.Range("G3").Value = "española"

# Selection moved to the whole column H.
$ws.Columns.Item(8).Select() | Out-Null
